$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 0.6864069264069264
$ws.Range("C4").Value = 0.1595858452183304
$ws.Range("D4").Value = 0.4
$ws.Range("E4").Value = 0.3571428571428571
$ws.Range("F4").Value = 0.7906976744186046
$ws.Range("G4").Value = 0.5
$ws.Range("I4").Value = 0.4166983114351535
$ws.Range("J4").Value = 0.738562091503268
$ws.Range("K4").Value = 31
$ws.Range("L4").Value = 3
$ws.Range("O4").Value = 0.8683333333333334
$ws.Range("P4").Value = 0.06311365409587162

# Row 5
$ws.Range("B5").Value = 0.7406060606060606
$ws.Range("C5").Value = 0.1070790924865831
$ws.Range("D5").Value = 0.4800000000000001
$ws.Range("E5").Value = 0.576923076923077
$ws.Range("F5").Value = 0.6976744186046512
$ws.Range("G5").Value = 0.375
$ws.Range("H5").Value = 0.6666666666666666
$ws.Range("I5").Value = 0.3809355902923154
$ws.Range("J5").Value = 0.7026143790849673
$ws.Range("K5").Value = 24
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 0.8922222222222222
$ws.Range("P5").Value = 0.05693834655697311

# Row 6
$ws.Range("B6").Value = 0.7311688311688311
$ws.Range("C6").Value = 0.05510603165703743
$ws.Range("I6").Value = 0.4372775372775373
$ws.Range("J6").Value = 0.7222222222222222
$ws.Range("O6").Value = 0.9122222222222222
$ws.Range("P6").Value = 0.04042978977480057

# Row 7
$ws.Range("B7").Value = 0.6303896103896104
$ws.Range("C7").Value = 0.1432146647838185
$ws.Range("D7").Value = 0.4285714285714285
$ws.Range("E7").Value = 0.3658536585365854
$ws.Range("F7").Value = 0.813953488372093
$ws.Range("G7").Value = 0.6
$ws.Range("I7").Value = 0.5066968130921619
$ws.Range("J7").Value = 0.7287581699346404
$ws.Range("K7").Value = 32
$ws.Range("L7").Value = 2
$ws.Range("O7").Value = 0.8300000000000001
$ws.Range("P7").Value = 0.1127435635019184

# Row 8
$ws.Range("B8").Value = 0.5944444444444443
$ws.Range("C8").Value = 0.2370081000855727
$ws.Range("I8").Value = 0.4096587390065651
$ws.Range("O8").Value = 0.8166666666666667
$ws.Range("P8").Value = 0.06411794687223779
